$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# "Supervision" block: mark item 22 (row 22, column D - "Note 0/1") as achieved (0 -> 1)
$ws.Range("D22").Value = 1

# Move the active selection, as recorded in the saved view state
$ws.Range("K19").Select()
